# Append a new data row (row 89) to the sheet, mirroring the existing
# Adafruit IO feed export rows. All values are stored as text, matching
# the source data (inline/shared strings), not as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 89

# Force text formatting on the new row first so that numeric-looking
# strings (e.g. "25") are kept as text instead of being auto-converted
# to numbers when the Value is assigned.
$ws.Range("A$row`:F$row").NumberFormat = "@"

$ws.Range("A$row").Value = "2024-09-25T18:06:40Z"
$ws.Range("B$row").Value = "temperature"
$ws.Range("C$row").Value = "25"
$ws.Range("D$row").Value = "N/A"
$ws.Range("E$row").Value = "N/A"
$ws.Range("F$row").Value = "N/A"
